$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings) ---
$ws.Range("C1").Value = "Integral"
$ws.Range("D1").Value = "Time"
$ws.Range("A13").Value = "Avg"

# --- Data updates (B2:D11) ---
$data = @(
    @(0.009288082218, 0.188861165865838, 0.003116144),
    @(0.009025169295, 0.200714787035528, 0.002247606),
    @(0.00818141143,  0.191453408375534, 0.001738233),
    @(0.0088677669,   0.189237485091271, 0.001496893),
    @(0.009042464821, 0.195915058494926, 0.001489197),
    @(0.007439458877, 0.188981557035694, 0.001500384),
    @(0.008710050355, 0.185580337153779, 0.001486881),
    @(0.008398276459, 0.178167083645935, 0.001612146),
    @(0.01083360915,  0.20395114276114,  0.001432878),
    @(0.01323649942,  0.184220097511442, 0.001849358)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
    $ws.Cells.Item($row, 4).Value = $data[$i][2]
}

# --- Formulas: row 13 average, row 14 sample stdev (not shared anymore) ---
$ws.Range("B13").Formula = "=AVERAGE(B2:B11)"
$ws.Range("C13").Formula = "=AVERAGE(C2:C11)"
$ws.Range("D13").Formula = "=AVERAGE(D2:D11)"

$ws.Range("B14").Formula = "=_xlfn.STDEV.S(B2:B11)"
$ws.Range("C14").Formula = "=_xlfn.STDEV.S(C2:C11)"
$ws.Range("D14").Formula = "=_xlfn.STDEV.S(D2:D11)"

# --- Remove extra sheet "Ark1" ---
$wb.Worksheets.Item("Ark1").Delete()

# --- Selection change on sheet1 ---
$ws.Range("D14").Select()
